$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text / summary value updates -------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 4

# --- Saturday (07/12/2025) section --------------------------------------
# Existing line item "Point 14 / POL-35-5" pricing goes to 0
$ws.Range("H16").Value = 0

# Insert a new line item row (shifts the old TOTAL row, and everything
# below it, down by one row) directly beneath the existing data row.
$ws.Rows.Item(17).Insert()

# Copy the formatting of the "alternate" data-row style (the grey banded
# row used for the 2nd+ line item of a day) onto the freshly inserted row
# so the new row matches the report's existing visual pattern.
$ws.Range("A24:H24").Copy()
$ws.Range("A17:H17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new line item: Point 14 / SVC-VA
$ws.Range("A17").Value = "Point 14"
$ws.Range("B17").Value = "SVC-VA"
$ws.Range("C17").Value = "Rem"
$ws.Range("D17").Value = "SVC-Virtual Asset Capitalization"
$ws.Range("E17").Value = "EA"
$ws.Range("F17").Value = 1
$ws.Range("H17").Value = 0

# Saturday TOTAL row (now row 18 after the insert above)
$ws.Range("H18").Value = 0

# --- Sunday (07/13/2025) section (rows shifted down by 1 to 21-25) ------
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
